$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 into the two new header cells, then set their text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I (I0) and J (IF)
$values = @(
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(4, 6),
    @(1, 4),
    @(1, 3),
    @(4, 6),
    @(3, 4),
    @(1, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
